$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Copy the header style from an existing header cell (A1) so the new
# header cells (AD1:AF1) pick up the same bold/centered/bordered format.
$ws.Range("A1").Copy()
$ws.Range("AD1:AF1").PasteSpecial(-4122)

$ws.Range("AD1").Value = "Wins"
$ws.Range("AE1").Value = "Losses"
$ws.Range("AF1").Value = "Ties"

For ($row = 2; $row -le 40; $row++) {
    $ws.Cells.Item($row, 30).Value = 73
    $ws.Cells.Item($row, 31).Value = 89
    $ws.Cells.Item($row, 32).Value = 0
}
